$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Find the row whose column B holds the image name that must be removed.
$target = "00063_batch2.jpg"
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row  # xlUp

$rowToDelete = -1
for ($r = 1; $r -le $lastRow; $r++) {
    if ($ws.Cells.Item($r, 2).Value() -eq $target) {
        $rowToDelete = $r
        break
    }
}

if ($rowToDelete -ge 1) {
    # Remove the entire row; remaining rows shift up automatically.
    $ws.Rows.Item($rowToDelete).Delete()

    # Column A holds a simple sequential index (0-based) starting on row 2.
    # Re-sequence it from the deletion point through the new last row so it
    # stays contiguous after the row shift.
    $newLastRow = $lastRow - 1
    for ($r = $rowToDelete; $r -le $newLastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 2
    }
}
